$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook window geometry (bookViews/workbookView) - best effort.
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Left = 420
$wb.Windows.Item(1).Top = 800
$wb.Windows.Item(1).Width = 43920
$wb.Windows.Item(1).Height = 16420

# ---------------------------------------------------------------------------
# Sheet "90-90-90" (sheet1): update row 6 (DRC) data + note, then selection.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("90-90-90")
$ws1.Activate()

$ws1.Range("B6").Value = 0.59
$ws1.Range("C6").Value = 0.74
$ws1.Range("D6").Value = 0.39
$ws1.Range("E6").Value = 0.49
$ws1.Range("F6").Value = 0.64
$ws1.Range("G6").Value = 0.07
$ws1.Range("H6").Value = 0.67
$ws1.Range("I6").Value = 0.83
$ws1.Range("J6").Value = 0.7
$ws1.Range("K6").Value = "max_error = 0.15, min_number = 1000"

# ---------------------------------------------------------------------------
# Sheet "2015" (sheet2): update row 6 data, then selection.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2015")
$ws2.Activate()

$ws2.Range("B6").Value = 393501
$ws2.Range("C6").Value = 165837
$ws2.Range("D6").Value = 98597
$ws2.Range("E6").Value = 97990
$ws2.Range("F6").Value = 36986
$ws2.Range("G6").Value = 348228
$ws2.Range("H6").Value = 128922
$ws2.Range("I6").Value = 75035
$ws2.Range("J6").Value = 74501
$ws2.Range("K6").Value = 6021
$ws2.Range("L6").Value = 444139
$ws2.Range("M6").Value = 200101
$ws2.Range("N6").Value = 120024
$ws2.Range("O6").Value = 119082
$ws2.Range("P6").Value = 74557

$ws2.Range("B6:P6").Select()

# ---------------------------------------------------------------------------
# Sheet "2020" (sheet3): update row 6 data, then selection.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("2020")
$ws3.Activate()

$ws3.Range("B6").Value = 353322
$ws3.Range("C6").Value = 206972
$ws3.Range("D6").Value = 153230
$ws3.Range("E6").Value = 152890
$ws3.Range("F6").Value = 61210
$ws3.Range("G6").Value = 305348
$ws3.Range("H6").Value = 166117
$ws3.Range("I6").Value = 116445
$ws3.Range("J6").Value = 115914
$ws3.Range("K6").Value = 9767
$ws3.Range("L6").Value = 401548
$ws3.Range("M6").Value = 245245
$ws3.Range("N6").Value = 190544
$ws3.Range("O6").Value = 190204
$ws3.Range("P6").Value = 123451

$ws3.Range("B6:P6").Select()

# ---------------------------------------------------------------------------
# Sheet "AIDS Deaths" (sheet4): update row 6 data, then selection.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("AIDS Deaths")
$ws4.Activate()

$ws4.Range("B6").Value = 17897
$ws4.Range("C6").Value = 19976
$ws4.Range("D6").Value = 21315
$ws4.Range("E6").Value = 21820
$ws4.Range("F6").Value = 22079
$ws4.Range("G6").Value = 14633
$ws4.Range("H6").Value = 15515
$ws4.Range("I6").Value = 16613
$ws4.Range("J6").Value = 17092
$ws4.Range("K6").Value = 17506
$ws4.Range("L6").Value = 21769
$ws4.Range("M6").Value = 25113
$ws4.Range("N6").Value = 26309
$ws4.Range("O6").Value = 26651
$ws4.Range("P6").Value = 26557

$ws4.Range("B6:P6").Select()

# ---------------------------------------------------------------------------
# Sheet "New Infections" (sheet5): update row 6 data, then selection.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("New Infections")
$ws5.Activate()

$ws5.Range("B6").Value = 15107
$ws5.Range("C6").Value = 15164
$ws5.Range("D6").Value = 15060
$ws5.Range("E6").Value = 14903
$ws5.Range("F6").Value = 14710
$ws5.Range("G6").Value = 9468
$ws5.Range("H6").Value = 9323
$ws5.Range("I6").Value = 9070
$ws5.Range("J6").Value = 8794
$ws5.Range("K6").Value = 8503
$ws5.Range("L6").Value = 20641
$ws5.Range("M6").Value = 20868
$ws5.Range("N6").Value = 20809
$ws5.Range("O6").Value = 20951
$ws5.Range("P6").Value = 21246

$ws5.Range("C31").Select()

# ---------------------------------------------------------------------------
# Re-activate the "90-90-90" sheet and restore its selection, so it is the
# tab that is active/selected when the workbook is reopened - matching the
# original tabSelected="1" on sheet1.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B34").Select()
